$d = $word.ActiveDocument

# The sentence "Index.php is genereal file, index.php have spl_autoload_register
# function. Autoload function includes our classes" was split across three runs
# (the middle one being just a single space). Re-run Find & Replace over the
# full sentence so Word rewrites it as a single, unified run while preserving
# the existing character formatting.
$text = "Index.php is genereal file, index.php have spl_autoload_register function. Autoload function includes our classes"

$d.Content.Find.Execute(
    $text,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $text,
    2
)
